$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11 ---
$ws.Range("A11").Value = 10
$ws.Range("H11").Value = "y"
$ws.Range("I11").Value = "Error creating bean with name 'adminController': Unsatisfied dependency expressed through field 'categoryDAO'"

# --- Row 12 ---
$ws.Range("A12").Value = 11
$ws.Range("H12").Value = "n"
$ws.Range("I12").Value = 'Sequence "HIBERNATE_SEQUENCE" not found; SQL statement: sol:made auto generate type as identity'

# --- Row 13 ---
$ws.Range("A13").Value = 12
$ws.Range("H13").Value = "y"
$ws.Range("I13").Value = 'org.springframework.dao.DataIntegrityViolationException: could not execute statement; SQL [n/a]; constraint ["FK_USERROLE_USER_ID: PUBLIC.USERROLE FOREIGN KEY(USER_ID) REFERENCES PUBLIC.USER(USER_ID) (0)"; SQL statement: sol: check if field names are proper'
$ws.Range("G13").Value = "20m"

# --- Row 14 ---
$ws.Range("A14").Value = 13
$ws.Range("H14").Value = "n"
$ws.Range("I14").Value = 'Sequence "HIBERNATE_SEQUENCE" not found; SQL statement: sol:made auto generate type as identity'

# --- Row 15 ---
$ws.Range("A15").Value = 14
$ws.Range("I15").Value = "java.io.StreamCorruptedException: invalid type code:"

# --- Row 16 ---
$ws.Range("A16").Value = 15
$ws.Range("I16").Value = "java.lang.ClassCastException: cannot assign instance of java.lang.StackTraceElement to field java.util.Collections`$UnmodifiableList.list of type java.util.List in instance of java.util.Collections`$UnmodifiableList"

# --- Dates for the new rows (reuse the existing date style/format from C2) ---
$ws.Range("C11:C16").Value = 42711
$ws.Range("C2").Copy()
$ws.Range("C11:C16").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Selection, matching the diff's new active cell ---
[void]$ws.Range("A10").Select()
